# Update recomputed NATMI ligand-receptor edge-weight statistics for Angptl3-Itgb3
# (new TPM-based expression values) as described in the commit "update scripts wuth new tpm".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: FAPs -> ECs
$ws.Range("I2").Value = 0.4458169960283037
$ws.Range("J2").Value = 0.4458169960283037
$ws.Range("M2").Value = 2.598166333333333
$ws.Range("N2").Value = 7.794499
$ws.Range("O2").Value = 0.3466013321552429
$ws.Range("P2").Value = 0.3466013321552429
$ws.Range("Q2").Value = 6.239352684296222
$ws.Range("R2").Value = 56.154174158666
$ws.Range("S2").Value = 0.1545207647208587
$ws.Range("T2").Value = 0.1545207647208587

# Row 3: FAPs -> FAPs
$ws.Range("I3").Value = 0.4458169960283037
$ws.Range("J3").Value = 0.4458169960283037
$ws.Range("M3").Value = 4.333403333333333
$ws.Range("O3").Value = 0.5780859172985858
$ws.Range("P3").Value = 0.5780859172985858
$ws.Range("R3").Value = 93.65785491013999
$ws.Range("S3").Value = 0.2577205270963219
$ws.Range("T3").Value = 0.2577205270963219

# Row 4: FAPs -> MuSCs
$ws.Range("I4").Value = 0.4458169960283037
$ws.Range("J4").Value = 0.4458169960283037
$ws.Range("M4").Value = 0.4692043333333333
$ws.Range("N4").Value = 1.407613
$ws.Range("O4").Value = 0.06259293136852516
$ws.Range("P4").Value = 0.06259293136852516
$ws.Range("Q4").Value = 1.126768243860222
$ws.Range("R4").Value = 10.140914194742
$ws.Range("S4").Value = 0.02790499263532167
$ws.Range("T4").Value = 0.02790499263532166

# Row 5: FAPs -> Resolving-Mac
$ws.Range("I5").Value = 0.4458169960283037
$ws.Range("J5").Value = 0.4458169960283037
$ws.Range("K5").Value = 2.0
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.09534933333333333
$ws.Range("N5").Value = 0.286048
$ws.Range("O5").Value = 0.01271981917764605
$ws.Range("P5").Value = 0.01271981917764604
$ws.Range("Q5").Value = 0.2289761480035555
$ws.Range("R5").Value = 2.060785332032
$ws.Range("S5").Value = 0.005670711575801369
$ws.Range("T5").Value = 0.005670711575801367

# Row 6: MuSCs -> ECs
$ws.Range("G6").Value = 1.399743666666667
$ws.Range("H6").Value = 4.199231
$ws.Range("I6").Value = 0.2598558798146963
$ws.Range("J6").Value = 0.2598558798146962
$ws.Range("M6").Value = 2.598166333333333
$ws.Range("N6").Value = 7.794499
$ws.Range("O6").Value = 0.3466013321552429
$ws.Range("P6").Value = 0.3466013321552429
$ws.Range("Q6").Value = 3.636766870029889
$ws.Range("R6").Value = 32.730901830269
$ws.Range("S6").Value = 0.09006639411214643
$ws.Range("T6").Value = 0.0900663941121464

# Row 7: MuSCs -> FAPs
$ws.Range("G7").Value = 1.399743666666667
$ws.Range("H7").Value = 4.199231
$ws.Range("I7").Value = 0.2598558798146963
$ws.Range("J7").Value = 0.2598558798146962
$ws.Range("M7").Value = 4.333403333333333
$ws.Range("O7").Value = 0.5780859172985858
$ws.Range("P7").Value = 0.5780859172985858
$ws.Range("Q7").Value = 6.065653870945555
$ws.Range("R7").Value = 54.59088483851
$ws.Range("S7").Value = 0.1502190246481098
$ws.Range("T7").Value = 0.1502190246481097

# Row 8: MuSCs -> MuSCs
$ws.Range("G8").Value = 1.399743666666667
$ws.Range("H8").Value = 4.199231
$ws.Range("I8").Value = 0.2598558798146963
$ws.Range("J8").Value = 0.2598558798146962
$ws.Range("M8").Value = 0.4692043333333333
$ws.Range("N8").Value = 1.407613
$ws.Range("O8").Value = 0.06259293136852516
$ws.Range("P8").Value = 0.06259293136852516
$ws.Range("Q8").Value = 0.656765793955889
$ws.Range("R8").Value = 5.910892145603
$ws.Range("S8").Value = 0.016265141250949
$ws.Range("T8").Value = 0.016265141250949

# Row 9: MuSCs -> Resolving-Mac
$ws.Range("G9").Value = 1.399743666666667
$ws.Range("H9").Value = 4.199231
$ws.Range("I9").Value = 0.2598558798146963
$ws.Range("J9").Value = 0.2598558798146962
$ws.Range("K9").Value = 2.0
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.09534933333333333
$ws.Range("N9").Value = 0.286048
$ws.Range("O9").Value = 0.01271981917764605
$ws.Range("P9").Value = 0.01271981917764604
$ws.Range("Q9").Value = 0.1334646254542222
$ws.Range("R9").Value = 1.201181629088
$ws.Range("S9").Value = 0.003305319803491059
$ws.Range("T9").Value = 0.003305319803491059

# Row 10: Resolving-Mac -> ECs
$ws.Range("G10").Value = 1.585427
$ws.Range("H10").Value = 4.756281
$ws.Range("I10").Value = 0.294327124157
$ws.Range("J10").Value = 0.294327124157
$ws.Range("M10").Value = 2.598166333333333
$ws.Range("N10").Value = 7.794499
$ws.Range("O10").Value = 0.3466013321552429
$ws.Range("P10").Value = 0.3466013321552429
$ws.Range("Q10").Value = 4.119203055357667
$ws.Range("R10").Value = 37.072827498219
$ws.Range("S10").Value = 0.1020141733222378
$ws.Range("T10").Value = 0.1020141733222378

# Row 11: Resolving-Mac -> FAPs
$ws.Range("G11").Value = 1.585427
$ws.Range("H11").Value = 4.756281
$ws.Range("I11").Value = 0.294327124157
$ws.Range("J11").Value = 0.294327124157
$ws.Range("M11").Value = 4.333403333333333
$ws.Range("O11").Value = 0.5780859172985858
$ws.Range("P11").Value = 0.5780859172985858
$ws.Range("Q11").Value = 6.870294646556665
$ws.Range("R11").Value = 61.83265181900999
$ws.Range("S11").Value = 0.1701463655541541
$ws.Range("T11").Value = 0.1701463655541541

# Row 12: Resolving-Mac -> MuSCs
$ws.Range("G12").Value = 1.585427
$ws.Range("H12").Value = 4.756281
$ws.Range("I12").Value = 0.294327124157
$ws.Range("J12").Value = 0.294327124157
$ws.Range("M12").Value = 0.4692043333333333
$ws.Range("N12").Value = 1.407613
$ws.Range("O12").Value = 0.06259293136852516
$ws.Range("P12").Value = 0.06259293136852516
$ws.Range("Q12").Value = 0.7438892185836666
$ws.Range("R12").Value = 6.695002967252999
$ws.Range("S12").Value = 0.01842279748225448
$ws.Range("T12").Value = 0.01842279748225448

# Row 13: Resolving-Mac -> Resolving-Mac
$ws.Range("G13").Value = 1.585427
$ws.Range("H13").Value = 4.756281
$ws.Range("I13").Value = 0.294327124157
$ws.Range("J13").Value = 0.294327124157
$ws.Range("K13").Value = 2.0
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.09534933333333333
$ws.Range("N13").Value = 0.286048
$ws.Range("O13").Value = 0.01271981917764605
$ws.Range("P13").Value = 0.01271981917764604
$ws.Range("Q13").Value = 0.1511694074986666
$ws.Range("R13").Value = 1.360524667488
$ws.Range("S13").Value = 0.003743787798353618
$ws.Range("T13").Value = 0.003743787798353616
